$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Tue Feb 11 20:20:53 EST 2025"
$ws.Range("B3").Value = "Tue Feb 11 20:21:06 EST 2025"
$ws.Range("B4").Value = "Tue Feb 11 20:21:18 EST 2025"
$ws.Range("B5").Value = "Tue Feb 11 20:21:30 EST 2025"
$ws.Range("B6").Value = "Tue Feb 11 20:21:42 EST 2025"
$ws.Range("B7").Value = "Tue Feb 11 20:21:54 EST 2025"
